# Auto-generated edit script: restores earlier ITSx default values
# and re-generated output numbers for mock1 and mock2 sheets.

$wb = $excel.ActiveWorkbook

# --- mock1 ---
$values_sheet1 = @{
    "C2" = 8864
    "D2" = 8865
    "E2" = 1726
    "C3" = 2283
    "D3" = 2283
    "I3" = 0.257558664259928
    "J3" = 0.257529610829103
    "I4" = 0.170803249097473
    "J4" = 0.170783981951495
    "K4" = 0.249710312862109
    "C5" = 1233
    "D5" = 1234
    "I5" = 0.139101985559567
    "J5" = 0.139199097574732
    "K5" = 0.198725376593279
    "I6" = 0.128497292418773
    "J6" = 0.128482797518331
    "K6" = 0.0926998841251448
    "C7" = 931
    "D7" = 931
    "E7" = 258
    "I7" = 0.105031588447653
    "J7" = 0.105019740552735
    "K7" = 0.149478563151796
    "C8" = 808
    "D8" = 808
    "I8" = 0.0911552346570397
    "J8" = 0.0911449520586576
    "K8" = 0.129200463499421
    "E9" = 199
    "I9" = 0.0581001805054152
    "J9" = 0.0580936266215454
    "K9" = 0.115295480880649
    "C10" = 158
    "D10" = 158
    "I10" = 0.0178249097472924
    "J10" = 0.0178228990411732
    "K10" = 0.0243337195828505
    "I11" = 0.0107175090252708
    "J11" = 0.0107163000564016
    "K11" = 0.0208574739281576
    "I12" = 0.0106046931407942
    "J12" = 0.0106034968979131
    "I13" = 0.0039485559566787
    "J13" = 0.00394811054709532
    "K13" = 0.00869061413673233
    "I14" = 0.00259476534296029
    "J14" = 0.00259447264523407
    "I15" = 0.00146660649819495
    "J15" = 0.00146644106034969
    "K15" = 0.00405561993047509
    "I16" = 0.00135379061371841
    "J16" = 0.00135363790186125
    "K16" = 0.00231749710312862
    "I17" = 0.00124097472924188
    "J17" = 0.00124083474337281
    "K17" = 0.00463499420625724
}
$ws = $wb.Worksheets.Item("mock1")
foreach ($addr in $values_sheet1.Keys) {
    $ws.Range($addr).Value = $values_sheet1[$addr]
}

# --- mock2 ---
$values_sheet2 = @{
    "C2" = 4638
    "D2" = 4630
    "E2" = 877
    "C3" = 1961
    "D3" = 1961
    "I3" = 0.422811556705476
    "J3" = 0.42354211663067
    "K3" = 0.673888255416192
    "C4" = 1494
    "D4" = 1494
    "I4" = 0.322121604139715
    "J4" = 0.32267818574514
    "C5" = 617
    "D5" = 617
    "I5" = 0.133031479085813
    "J5" = 0.133261339092873
    "K5" = 0.182440136830103
    "I6" = 0.0379473911168607
    "J6" = 0.0380129589632829
    "K6" = 0.0592930444697834
    "E7" = 44
    "I7" = 0.0351444588184562
    "J7" = 0.0352051835853132
    "K7" = 0.0501710376282782
    "I8" = 0.0288917636912462
    "J8" = 0.0289416846652268
    "K8" = 0.0182440136830103
    "I9" = 0.0103492884864166
    "J9" = 0.0103671706263499
    "K9" = 0.00456100342075257
    "I10" = 0.00495903406640793
    "J10" = 0.00496760259179266
    "K10" = 0.00798175598631699
    "I11" = 0.00301854247520483
    "J11" = 0.00302375809935205
    "K11" = 0.00342075256556442
    "I12" = 0.00172488141440276
}
$ws = $wb.Worksheets.Item("mock2")
foreach ($addr in $values_sheet2.Keys) {
    $ws.Range($addr).Value = $values_sheet2[$addr]
}
